# Update "RELAÇÃO DE FISCAIS DE CONTRATOS VIGENTES.xlsx" - Plan2 sheet
# - Fill in two missing "Valor Pago" values (I6, I16)
# - Add the "Saldo" formula (H-I) down column J for rows 6:60
# - Leave the selection on J16, matching the saved workbook state

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan2")

# New data entered by the author
$ws.Range("I6").Value = 167567
$ws.Range("I16").Value = 1453410.42

# New "Saldo" formulas: J6 is a standalone formula, J7:J60 is one shared-formula block
$ws.Range("J6").Formula = "=H6-I6"
$ws.Range("J7:J60").Formula = "=H7-I7"

# Final UI state: active sheet Plan2, selection on J16
$ws.Activate()
$ws.Range("J16").Select()
